$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Tags term accession number / term source ref (rows 13/14, columns C = NCIT, for "Proteomics")
$ws.Range("C13").Value = "http://purl.obolibrary.org/obo/NCIT_C20085"
$ws.Range("C14").Value = "NCIT"

# ER term accession number / term source ref (row 9 / 10, column B = DPBO, for "PRIDE")
$ws.Range("B9").Value = "http://purl.obolibrary.org/obo/DPBO_1000098"
$ws.Range("B10").Value = "DPBO"

# Tags term accession number / term source ref (rows 13/14, columns B = DPBO, for "PRIDE")
$ws.Range("B13").Value = "http://purl.obolibrary.org/obo/DPBO_1000098"
$ws.Range("B14").Value = "DPBO"

# Select the SwateTemplateMetadata sheet and set the active cell, matching
# the diff's updated tab selection / active cell.
$ws.Activate()
$ws.Range("B11").Select()
